# Rename the inline logo pictures that live in the headers/footers.
#
# Pearson logo (appears in both the default and first-page footers) goes
# from "image2.png" -> "image1.png".
# BTEC logo (appears in the first-page header) goes from
# "image1.jpg" -> "image2.jpg".
#
# We identify each picture by its (stable) alt text / description rather
# than by a hard-coded header/footer index, so the script keeps working
# even if the section/header/footer enumeration order differs.

$d = $word.ActiveDocument

function Rename-LogoShapes($range) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $descr = $shape.AlternativeText

        if ($descr -like "*PearsonLogo.png") {
            $shape.Name = "image1.png"
        }
        elseif ($descr -eq "BTec_Logo-Orange") {
            $shape.Name = "image2.jpg"
        }
    }
}

for ($secIndex = 1; $secIndex -le $d.Sections.Count; $secIndex++) {
    $section = $d.Sections.Item($secIndex)

    for ($hIndex = 1; $hIndex -le $section.Headers.Count; $hIndex++) {
        $header = $section.Headers.Item($hIndex)
        if ($header.Exists) {
            Rename-LogoShapes($header.Range)
        }
    }

    for ($fIndex = 1; $fIndex -le $section.Footers.Count; $fIndex++) {
        $footer = $section.Footers.Item($fIndex)
        if ($footer.Exists) {
            Rename-LogoShapes($footer.Range)
        }
    }
}
